$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the old row 12 ("Programa resumido:"), shifting
# rows 12-20 down to 15-23 (and carrying their per-row heights/styles with
# them, same as Excel's native Insert behaviour).
$ws.Rows("12:14").Insert()

# New row 12: "Docentes responsáveis:" label only (column A).
$ws.Range("A12").Value = "Docentes responsáveis:"

# New row 13: first responsible professor, duplicated in B and C.
$ws.Range("B13:C13").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# New row 14: second responsible professor, duplicated in B and C.
$ws.Range("B14:C14").Value = "7797767 - Viktor Pastoukhov"

# The Insert() call copies formatting from the row above into the blank
# B12/C12/A13/A14 cells; drop those so no empty-but-styled cells remain.
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
